$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FATKP"
$ws.Range("B2").Value = 518219
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 3

$ws.Range("A1").Select()
